$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.2645
$ws.Range("E3").Value = 12.74229999999999
$ws.Range("E5").Value = 12.7185
$ws.Range("B9").Value = 8.633300000000004
$ws.Range("E11").Value = 13.25719999999999
$ws.Range("E12").Value = 13.04649999999999
$ws.Range("B13").Value = 5.698200000000003
$ws.Range("B16").Value = 9.084000000000009
$ws.Range("B18").Value = 5.297000000000002
$ws.Range("B20").Value = 5.425599999999998
$ws.Range("E21").Value = 12.80379999999999
